$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 41
$ws.Range("B17").Value = "changed readme file"
$ws.Range("C17").Value = "riya-morankar"
$ws.Range("D17").Value = "N/A"
$ws.Range("E17").Value = "edit1 to main"

# Force the date-looking string to stay as plain text (matches other rows),
# instead of being auto-converted into a date serial value.
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "2025-06-18"
$ws.Range("F17").Style = "Normal"
